$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.5284097792228408
$ws.Cells.Item(2, 3).Value = 0.24424712609407
$ws.Cells.Item(2, 5).Value = 0.1286485341705763
$ws.Cells.Item(2, 6).Value = 0.4443680307746263
$ws.Cells.Item(2, 7).Value = 0.5540818712052058
$ws.Cells.Item(2, 8).Value = 0.7071370714610552
$ws.Cells.Item(2, 11).Value = 0.2999255099659308
$ws.Cells.Item(2, 12).Value = 0.1851602229599507
$ws.Cells.Item(2, 15).Value = 2.496253976096042

$ws.Cells.Item(3, 2).Value = 0.4822756424105989
$ws.Cells.Item(3, 3).Value = 0.2460792773768397
$ws.Cells.Item(3, 5).Value = 0.128436951533601
$ws.Cells.Item(3, 6).Value = 0.387822817061874
$ws.Cells.Item(3, 7).Value = 0.5619452859039384
$ws.Cells.Item(3, 8).Value = 0.7146193146290969
$ws.Cells.Item(3, 11).Value = 0.2619185531828805
$ws.Cells.Item(3, 12).Value = 0.1777248589608149
$ws.Cells.Item(3, 15).Value = 2.528357410731232

$ws.Cells.Item(4, 2).Value = 0.453990798660044
$ws.Cells.Item(4, 3).Value = 0.247282015458385
$ws.Cells.Item(4, 5).Value = 0.128379939886532
$ws.Cells.Item(4, 6).Value = 0.3531389305168915
$ws.Cells.Item(4, 7).Value = 0.5671717665763438
$ws.Cells.Item(4, 8).Value = 0.7195230603105287
$ws.Cells.Item(4, 11).Value = 0.238498681436198
$ws.Cells.Item(4, 12).Value = 0.1732488698059171
$ws.Cells.Item(4, 15).Value = 2.549553848343699

$ws.Cells.Item(5, 2).Value = 0.4424758190413627
$ws.Cells.Item(5, 3).Value = 0.2477917427569416
$ws.Cells.Item(5, 5).Value = 0.1283750619643591
$ws.Cells.Item(5, 6).Value = 0.3390132514313251
$ws.Cells.Item(5, 7).Value = 0.5694016687474885
$ws.Cells.Item(5, 8).Value = 0.7215992674799345
$ws.Cells.Item(5, 11).Value = 0.2289345714353317
$ws.Cells.Item(5, 12).Value = 0.1714474437192024
$ws.Cells.Item(5, 15).Value = 2.558564810918128

$ws.Cells.Item(6, 2).Value = 0.4405644740710954
$ws.Cells.Item(6, 3).Value = 0.2478775676903346
$ws.Cells.Item(6, 5).Value = 0.1283753612555145
$ws.Cells.Item(6, 6).Value = 0.336668177824194
$ws.Cells.Item(6, 7).Value = 0.5697779835720951
$ws.Cells.Item(6, 8).Value = 0.7219487255536698
$ws.Cells.Item(6, 11).Value = 0.22734525035996
$ws.Cells.Item(6, 12).Value = 0.1711496849631686
$ws.Cells.Item(6, 15).Value = 2.560083613695696

$ws.Cells.Item(7, 2).Value = 0.4538354566564635
$ws.Cells.Item(7, 3).Value = 0.2472888103972117
$ws.Cells.Item(7, 5).Value = 0.1283797997522704
$ws.Cells.Item(7, 6).Value = 0.3529483938344953
$ws.Cells.Item(7, 7).Value = 0.5672014347611274
$ws.Cells.Item(7, 8).Value = 0.7195507453422465
$ws.Cells.Item(7, 11).Value = 0.2383697779092842
$ws.Cells.Item(7, 12).Value = 0.1732244835930459
$ws.Cells.Item(7, 15).Value = 2.549673862255318

$ws.Cells.Item(8, 2).Value = 0.5124945605967071
$ws.Cells.Item(8, 3).Value = 0.2448627311471299
$ws.Cells.Item(8, 5).Value = 0.1285604634304072
$ws.Cells.Item(8, 6).Value = 0.4248636149813336
$ws.Cells.Item(8, 7).Value = 0.5567104533690852
$ws.Cells.Item(8, 8).Value = 0.7096527092133869
$ws.Cells.Item(8, 11).Value = 0.2868384151776411
$ws.Cells.Item(8, 12).Value = 0.1825780254420835
$ws.Cells.Item(8, 15).Value = 2.507014998922187

$ws.Cells.Item(9, 2).Value = 0.6278247408028221
$ws.Cells.Item(9, 3).Value = 0.2407205966691564
$ws.Cells.Item(9, 5).Value = 0.1294923505423746
$ws.Cells.Item(9, 6).Value = 0.5661985755041457
$ws.Cells.Item(9, 7).Value = 0.5393020571075411
$ws.Cells.Item(9, 8).Value = 0.692697666592899
$ws.Cells.Item(9, 11).Value = 0.3811985982644046
$ws.Cells.Item(9, 12).Value = 0.2016256849989588
$ws.Cells.Item(9, 15).Value = 2.435147365130376

$ws.Cells.Item(10, 2).Value = 0.7127073771032428
$ws.Cells.Item(10, 3).Value = 0.2380500310172593
$ws.Cells.Item(10, 5).Value = 0.1305282992863326
$ws.Cells.Item(10, 6).Value = 0.6702781546542269
$ws.Cells.Item(10, 7).Value = 0.5284465344367177
$ws.Cells.Item(10, 8).Value = 0.6817348783814765
$ws.Cells.Item(10, 11).Value = 0.4500800247942607
$ws.Cells.Item(10, 12).Value = 0.2160468241055753
$ws.Cells.Item(10, 15).Value = 2.389537756991942

$ws.Cells.Item(11, 2).Value = 0.7513486486080296
$ws.Cells.Item(11, 3).Value = 0.2369155149208027
$ws.Cells.Item(11, 5).Value = 0.1310756986417871
$ws.Cells.Item(11, 6).Value = 0.7176906081379002
$ws.Cells.Item(11, 7).Value = 0.5239292575672181
$ws.Cells.Item(11, 8).Value = 0.6770715088859944
$ws.Cells.Item(11, 11).Value = 0.4813141764478246
$ws.Cells.Item(11, 12).Value = 0.2226994702938754
$ws.Cells.Item(11, 15).Value = 2.370351697121762

$ws.Cells.Item(12, 2).Value = 0.765984347166409
$ws.Cells.Item(12, 3).Value = 0.2364974154299162
$ws.Cells.Item(12, 5).Value = 0.1312939142844911
$ws.Cells.Item(12, 6).Value = 0.7356546913071611
$ws.Cells.Item(12, 7).Value = 0.5222793150505964
$ws.Cells.Item(12, 8).Value = 0.6753521210571662
$ws.Cells.Item(12, 11).Value = 0.4931267270023625
$ws.Cells.Item(12, 12).Value = 0.2252318544717298
$ws.Cells.Item(12, 15).Value = 2.363311194702419

$ws.Cells.Item(13, 2).Value = 0.7628321646024858
$ws.Cells.Item(13, 3).Value = 0.2365869489667958
$ws.Cells.Item(13, 5).Value = 0.131246431989247
$ws.Cells.Item(13, 6).Value = 0.7317853510981394
$ws.Cells.Item(13, 7).Value = 0.5226319602775504
$ws.Cells.Item(13, 8).Value = 0.6757203527991535
$ws.Cells.Item(13, 11).Value = 0.4905833677928513
$ws.Cells.Item(13, 12).Value = 0.224685876264445
$ws.Cells.Item(13, 15).Value = 2.364817489289393

$ws.Cells.Item(14, 2).Value = 0.7525526796043494
$ws.Cells.Item(14, 3).Value = 0.2368808870033021
$ws.Cells.Item(14, 5).Value = 0.1310934324973978
$ws.Cells.Item(14, 6).Value = 0.7191683204515869
$ws.Cells.Item(14, 7).Value = 0.5237922995622526
$ws.Cells.Item(14, 8).Value = 0.6769291213758564
$ws.Cells.Item(14, 11).Value = 0.4822863096950982
$ws.Cells.Item(14, 12).Value = 0.2229075477501681
$ws.Cells.Item(14, 15).Value = 2.369767963260358

$ws.Cells.Item(15, 2).Value = 0.7462565736511806
$ws.Cells.Item(15, 3).Value = 0.2370624313673844
$ws.Cells.Item(15, 5).Value = 0.131001138314808
$ws.Cells.Item(15, 6).Value = 0.7114413442032514
$ws.Cells.Item(15, 7).Value = 0.5245109431142936
$ws.Cells.Item(15, 8).Value = 0.6776755861088901
$ws.Cells.Item(15, 11).Value = 0.4772021285971277
$ws.Cells.Item(15, 12).Value = 0.2218199820721907
$ws.Cells.Item(15, 15).Value = 2.372829559425838

$ws.Cells.Item(16, 2).Value = 0.7101825601860412
$ws.Cells.Item(16, 3).Value = 0.2381257878826091
$ws.Cells.Item(16, 5).Value = 0.130494055761627
$ws.Cells.Item(16, 6).Value = 0.6671810134426437
$ws.Cells.Item(16, 7).Value = 0.5287502300168683
$ws.Cells.Item(16, 8).Value = 0.6820461522846912
$ws.Cells.Item(16, 11).Value = 0.4480367190402319
$ws.Cells.Item(16, 12).Value = 0.2156139074678833
$ws.Cells.Item(16, 15).Value = 2.39082306195003

$ws.Cells.Item(17, 2).Value = 0.6880587832855269
$ws.Cells.Item(17, 3).Value = 0.2387986738025134
$ws.Cells.Item(17, 5).Value = 0.1302024632827816
$ws.Cells.Item(17, 6).Value = 0.6400460337125793
$ws.Cells.Item(17, 7).Value = 0.5314587970659517
$ws.Cells.Item(17, 8).Value = 0.6848102405605374
$ws.Cells.Item(17, 11).Value = 0.4301184728015528
$ws.Cells.Item(17, 12).Value = 0.2118302651122264
$ws.Cells.Item(17, 15).Value = 2.402261718359185

$ws.Cells.Item(18, 2).Value = 0.6753364361325396
$ws.Cells.Item(18, 3).Value = 0.2391932638290655
$ws.Cells.Item(18, 5).Value = 0.1300419150644743
$ws.Cells.Item(18, 6).Value = 0.6244449056556647
$ws.Cells.Item(18, 7).Value = 0.5330563083680744
$ws.Cells.Item(18, 8).Value = 0.6864305359758731
$ws.Cells.Item(18, 11).Value = 0.4198029681089963
$ws.Cells.Item(18, 12).Value = 0.2096627179232797
$ws.Cells.Item(18, 15).Value = 2.40898793815083

$ws.Cells.Item(19, 2).Value = 0.6710293530423712
$ws.Cells.Item(19, 3).Value = 0.2393281654119193
$ws.Cells.Item(19, 5).Value = 0.1299887879731934
$ws.Cells.Item(19, 6).Value = 0.619163680173358
$ws.Cells.Item(19, 7).Value = 0.5336039985397676
$ws.Cells.Item(19, 8).Value = 0.6869843729362302
$ws.Cells.Item(19, 11).Value = 0.4163087222275692
$ws.Cells.Item(19, 12).Value = 0.2089303220823382
$ws.Cells.Item(19, 15).Value = 2.411290562157618

$ws.Cells.Item(20, 2).Value = 0.6904136279042632
$ws.Cells.Item(20, 3).Value = 0.2387262613750778
$ws.Cells.Item(20, 5).Value = 0.1302327620884505
$ws.Cells.Item(20, 6).Value = 0.642933953830422
$ws.Cells.Item(20, 7).Value = 0.5311663645901987
$ws.Cells.Item(20, 8).Value = 0.6845128457288396
$ws.Cells.Item(20, 11).Value = 0.4320268790970943
$ws.Cells.Item(20, 12).Value = 0.2122321403213761
$ws.Cells.Item(20, 15).Value = 2.401028837487033

$ws.Cells.Item(21, 2).Value = 0.7555719368936025
$ws.Cells.Item(21, 3).Value = 0.2367942379845474
$ws.Cells.Item(21, 5).Value = 0.1311380757638858
$ws.Cells.Item(21, 6).Value = 0.7228739723491628
$ws.Cells.Item(21, 7).Value = 0.5234498327564836
$ws.Cells.Item(21, 8).Value = 0.67657281405981
$ws.Cells.Item(21, 11).Value = 0.4847237734829264
$ws.Cells.Item(21, 12).Value = 0.2234295293933144
$ws.Cells.Item(21, 15).Value = 2.368307785067842

$ws.Cells.Item(22, 2).Value = 0.7981742396340223
$ws.Cells.Item(22, 3).Value = 0.2355986635659377
$ws.Cells.Item(22, 5).Value = 0.1317934335929074
$ws.Cells.Item(22, 6).Value = 0.7751780083420101
$ws.Cells.Item(22, 7).Value = 0.5187602153741224
$ws.Cells.Item(22, 8).Value = 0.6716547452439698
$ws.Cells.Item(22, 11).Value = 0.5190756354554367
$ws.Cells.Item(22, 12).Value = 0.2308243899738756
$ws.Cells.Item(22, 15).Value = 2.348233343811387

$ws.Cells.Item(23, 2).Value = 0.7754352771011668
$ws.Cells.Item(23, 3).Value = 0.236230634717046
$ws.Cells.Item(23, 5).Value = 0.1314378369755822
$ws.Cells.Item(23, 6).Value = 0.7472568307830727
$ws.Cells.Item(23, 7).Value = 0.521230760388292
$ws.Cells.Item(23, 8).Value = 0.6742548003375148
$ws.Cells.Item(23, 11).Value = 0.5007497427102976
$ws.Cells.Item(23, 12).Value = 0.226870632042278
$ws.Cells.Item(23, 15).Value = 2.358827449519552

$ws.Cells.Item(24, 2).Value = 0.6893490122189974
$ws.Cells.Item(24, 3).Value = 0.2387589749242593
$ws.Cells.Item(24, 5).Value = 0.1302190419057325
$ws.Cells.Item(24, 6).Value = 0.6416283278902171
$ws.Cells.Item(24, 7).Value = 0.5312984477280978
$ws.Cells.Item(24, 8).Value = 0.6846472007947
$ws.Cells.Item(24, 11).Value = 0.4311641323714923
$ws.Cells.Item(24, 12).Value = 0.2120504284872737
$ws.Cells.Item(24, 15).Value = 2.401585755758717

$ws.Cells.Item(25, 2).Value = 0.5965962172230661
$ws.Cells.Item(25, 3).Value = 0.2417755313379963
$ws.Cells.Item(25, 5).Value = 0.1291784991568896
$ws.Cells.Item(25, 6).Value = 0.5279251897347166
$ws.Cells.Item(25, 7).Value = 0.5436721667612048
$ws.Cells.Item(25, 8).Value = 0.6970218962581498
$ws.Cells.Item(25, 11).Value = 0.3557479206367873
$ws.Cells.Item(25, 12).Value = 0.1963976012061437
$ws.Cells.Item(25, 15).Value = 2.453327007466541
